# corregida gramatica del informe
# Adds a new pair of columns (K: "Tabú", L: "Busqueda Dispersa") to the four
# statistics tables on Hoja1, plus a fifth summary/header row (97) further
# down the sheet, and updates the sheet view (zoom / scroll / selection).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# ---------------------------------------------------------------------------
# Helper: write a header cell (column title) - reuses the same underline
# style as the rest of the header row.
# ---------------------------------------------------------------------------
function Set-HeaderCell($addr, $text) {
    $cell = $ws.Range($addr)
    $cell.Value = $text
    $cell.Font.Underline = $true
}

# ---------------------------------------------------------------------------
# Helper: write a plain numeric data cell, optionally reproducing the
# underline style used by the first block of data rows (2-11).
# ---------------------------------------------------------------------------
function Set-DataCell($addr, $value, [bool]$styled) {
    $cell = $ws.Range($addr)
    $cell.Value = $value
    if ($styled) {
        $cell.Font.Underline = $true
    }
}

# ---------------------------------------------------------------------------
# Table 1 (rows 1-11) - every cell in this block carries the header style.
# ---------------------------------------------------------------------------
Set-HeaderCell "K1" "Tabú"
Set-HeaderCell "L1" "Busqueda Dispersa"

$table1 = @{
    2  = @(4,3)
    3  = @(4,3)
    4  = @(4,3)
    5  = @(4,3)
    6  = @(4,3)
    7  = @(4,3)
    8  = @(4,3)
    9  = @(4,3)
    10 = @(4,3)
    11 = @(4,3)
}
foreach ($r in $table1.Keys) {
    $vals = $table1[$r]
    Set-DataCell "K$r" $vals[0] $true
    Set-DataCell "L$r" $vals[1] $true
}

# ---------------------------------------------------------------------------
# Table 2 (rows 19-30) - only the header row (19) is styled.
# ---------------------------------------------------------------------------
Set-HeaderCell "K19" "Tabú"
Set-HeaderCell "L19" "Busqueda Dispersa"

$table2 = @{
    20 = @(10,9)
    21 = @(10,9)
    22 = @(10,9)
    23 = @(11,9)
    24 = @(10,9)
    25 = @(10,9)
    26 = @(10,9)
    27 = @(10,9)
    28 = @(12,9)
    29 = @(10,9)
    30 = @(10,9)
}
foreach ($r in $table2.Keys) {
    $vals = $table2[$r]
    Set-DataCell "K$r" $vals[0] $false
    Set-DataCell "L$r" $vals[1] $false
}

# ---------------------------------------------------------------------------
# Table 3 (rows 40-51) - only the header row (40) gains new cells; the
# remaining data rows (41-51) simply widen their "spans" automatically.
# ---------------------------------------------------------------------------
Set-HeaderCell "K40" "Tabú"
Set-HeaderCell "L40" "Busqueda Dispersa"

# ---------------------------------------------------------------------------
# Table 4 (rows 60-71) - same pattern as table 3.
# ---------------------------------------------------------------------------
Set-HeaderCell "K60" "Tabú"
Set-HeaderCell "L60" "Busqueda Dispersa"

# ---------------------------------------------------------------------------
# Table 5 (rows 80-91) - same pattern again.
# ---------------------------------------------------------------------------
Set-HeaderCell "K80" "Tabú"
Set-HeaderCell "L80" "Busqueda Dispersa"

# ---------------------------------------------------------------------------
# New header/summary row 97 (starts at column B, not A).
# ---------------------------------------------------------------------------
$row97Headers = @{
    "B97" = "Bap: No mejora"
    "C97" = "Bra: nveces"
    "D97" = "Bra: nomejora"
    "E97" = "Busqueda Local"
    "F97" = "Grasp"
    "G97" = "Bam"
    "H97" = "Recocido Simulado"
    "I97" = "VND"
    "J97" = "BVNS"
    "K97" = "Tabú"
    "L97" = "Busqueda Dispersa"
}
foreach ($addr in @("B97","C97","D97","E97","F97","G97","H97","I97","J97","K97","L97")) {
    Set-HeaderCell $addr $row97Headers[$addr]
}

# ---------------------------------------------------------------------------
# New column width for column L (best-fit-like width for the longer header).
# ---------------------------------------------------------------------------
$ws.Columns.Item(12).ColumnWidth = 16.5

# ---------------------------------------------------------------------------
# Sheet view: scroll position, zoom and selection.
# ---------------------------------------------------------------------------
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 28
$win.ScrollColumn = 8
$win.Zoom = 85
$ws.Range("O20:O21").Select()
